$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 51: C51 was empty, now gets value 2
$ws.Range("C51").Value = 2

# Row 52 needs a new observer record (id 51, code "CL", coldef type 2).
# The row currently has no fill style on column B ("B52" has style index 1
# instead of the index 2 used by the surrounding populated rows). Duplicating
# row 51 (which already carries the correct per-column styles) into row 52's
# position and then overwriting the values reproduces that styling exactly,
# rather than leaving the default/unstyled look.
$ws.Rows(51).Copy()
$ws.Rows(52).Insert(-4121)   # xlShiftDown
# The row that got pushed down to 53 is just a duplicate of the (empty,
# identically styled) row that used to be at 52, so removing it restores the
# original row count/layout losslessly.
$ws.Rows(53).Delete()
# Row insertion resets the custom row height; restore it to match its
# neighbours.
$ws.Rows(52).RowHeight = 23.25

$ws.Range("A52").Value = 51
$ws.Range("B52").Value = "CL"
$ws.Range("C52").Value = 2

# Rows 53-55: only column A (the running index) changes.
$ws.Range("A53").Value = 52
$ws.Range("A54").Value = 53
$ws.Range("A55").Value = 54
